# Weekly update: a new price-survey record (dated 2023-10-04) is inserted
# as row 31 ("Haba" / Terminal Hortofrutícola Agro Chillán), pushing the
# previously-existing rows 31-85 down to rows 32-86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 31 - shifts rows 31..85 down to 32..86
# (mirrors Excel's normal "Insert Sheet Rows" behaviour, carrying the
# neighbouring row's number formatting onto the new cells).
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record's data.
$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(31, 3).Value = "Ñuble"
$ws.Cells.Item(31, 4).Value = 45203
$ws.Cells.Item(31, 5).Value = 16
$ws.Cells.Item(31, 6).Value = 100112026
$ws.Cells.Item(31, 7).Value = "Haba"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 120
$ws.Cells.Item(31, 11).Value = 14000
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = 14667
$ws.Cells.Item(31, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(31, 16).Value = 587
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
